$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$win = $excel.ActiveWindow
$pane1 = $win.Panes.Item(1)
$pane1.Activate()
$ws.Range("E1").Select()
$pane2 = $win.Panes.Item(2)
$pane2.Activate()
$ws.Range("I7").Select()
